# Requirements Trace Matrix - various fixes, comments, and refactoring
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Status (column D) values that were previously blank ---
$ws.Range("D24").Value = "P"
$ws.Range("D44").Value = "D"
$ws.Range("D45").Value = "D"
$ws.Range("D52").Value = "D"
$ws.Range("D53").Value = "D"
$ws.Range("D79").Value = "D"
$ws.Range("D82").Value = "D"
$ws.Range("D83").Value = "P"
$ws.Range("D85").Value = "D"
$ws.Range("D107").Value = "D"
$ws.Range("D108").Value = "D"
$ws.Range("D109").Value = "D"
$ws.Range("D145").Value = "D"

# --- Hide rows that are now considered "done" and folded into the filter view ---
$hiddenRows = @(21,22,23,24,25,40,41,42,43,44,45,52,53,60,61,77,78,79,80,81,82,107,108,109,140,141,142,143,144,145)
foreach ($r in $hiddenRows) {
    $ws.Rows.Item($r).Hidden = $true
}

# Row 43 also got manually resized
$ws.Rows.Item(43).RowHeight = 30

# --- Insert a new "Partial:" summary row above the existing "Left:" row ---
$ws.Rows.Item(181).Insert()
$ws.Range("F181").Value = "Partial:"
$ws.Range("G181").Formula = "=COUNTIF(D:D, ""P"")"

# Update the "Left:" row's formula (now shifted to row 182) to subtract the new Partial count too
$ws.Range("G182").Formula = "= 178 - G180 -G181"

# --- Change the table's autofilter from "blank status" to "Class = ModelTransformationTag" ---
$tbl = $ws.ListObjects.Item(1)
$tbl.Range.AutoFilter(1, @("requirements_class_ModelTransformationTag"), 7)

# --- Adjust the frozen pane scroll position and active selection ---
$excel.ActiveWindow.ScrollRow = 41
$ws.Range("F185").Select()
